$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A to hold the serial number (S/No.)
$ws.Columns.Item(1).Insert()

# Headers
$ws.Range("A1").Value = "S/No."
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Price"
$ws.Range("D1").Value = "Quantity"
$ws.Range("E1").Value = "Category"

# Row 2 - Rice
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Rice"
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 100
$ws.Range("E2").Value = "Groceries"

# Row 3 - Beer
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Beer"
$ws.Range("C3").Value = 50
$ws.Range("D3").Value = 100
$ws.Range("E3").Value = "Drinks"

# Row 4 - Soap (replaces Shirts)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Soap"
$ws.Range("C4").Value = 80
$ws.Range("D4").Value = 70
$ws.Range("E4").Value = "Toiletries"

$ws.Range("F4").Select()
